$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update ledStatus for project Z005 rows (E3 and E6) from "Concluido" to numeric 2
$ws.Range("E3").Value = 2
$ws.Range("E6").Value = 2
